$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '21.738.65'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').Value = '1.539.42'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '289.90'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3892'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +3.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3187'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.92'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07204'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.059'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -5.95%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.645'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.64'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.615'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.52%  '
$ws.Range('D16').Value = '1.539.06'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001108'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06587'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.28'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.12%  '
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('E21').Value = '  -4.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.40'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.88'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -5.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.391'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.52%  '
$ws.Range('D25').Value = '21.736.02'
$ws.Range('E25').Value = '  -1.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.375'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -5.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '146.63'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.39'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.842'
$ws.Range('D29').ClearFormats()
$ws.Range('D30').Value = '1.717.84'
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.58'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9783'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -12.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.934'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08205'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.831'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -4.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06089'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.48%  '
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.479'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -9.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02205'
$ws.Range('D39').ClearFormats()
$ws.Range('E40').Value = '  -3.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.190'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.06%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.69'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5754'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.07'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.746'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5512'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.07%  '
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '117.11'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.872'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06720'
$ws.Range('D51').ClearFormats()
